$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.00379999999999
$ws.Range("D4").Value = -8.279500000000002

$ws.Range("C7").Value = -13.4302

$ws.Range("D12").Value = -5.935399999999998

$ws.Range("C16").Value = -14.181

$ws.Range("D18").Value = -8.758100000000001

$ws.Range("D19").Value = -8.394999999999992

$ws.Range("D20").Value = -8.451899999999991

$ws.Range("C28").Value = -13.628

$ws.Range("C29").Value = -11.6097

$ws.Range("D31").Value = -7.333199999999997

$ws.Range("C32").Value = -12.48889999999999

$ws.Range("C40").Value = -12.5757
$ws.Range("D40").Value = -8.405799999999994

$ws.Range("D42").Value = -8.563699999999997

$ws.Range("D47").Value = -7.424099999999999

$ws.Range("D48").Value = -7.196399999999996

$ws.Range("C52").Value = -11.2898

$ws.Range("C57").Value = -14.3094

$ws.Range("D63").Value = -6.596699999999998

$ws.Range("D64").Value = -7.015899999999993

$ws.Range("C66").Value = -11.4096

$ws.Range("D76").Value = -7.7175

$ws.Range("D81").Value = -7.483399999999996

$ws.Range("D89").Value = -8.211800000000004

$ws.Range("D94").Value = -6.053500000000001

$ws.Range("C100").Value = -12.0818
